# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
# Mirrors editing individual cells in the "cryptos" worksheet as a human would in Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.955.23"
$ws.Range("E2").Value = "  +1.21%  "

$ws.Range("D3").Value = "1.759.43"
$ws.Range("E3").Value = "  -0.94%  "

$ws.Range("E4").Value = "  -0.91%  "

$ws.Range("D5").Value = "'336.05"
$ws.Range("E5").Value = "  -0.53%  "

$ws.Range("E6").Value = "  -0.87%  "

$ws.Range("D7").Value = "'0.3842"
$ws.Range("E7").Value = "  -1.23%  "

$ws.Range("D8").Value = "'0.3390"
$ws.Range("E8").Value = "  -1.44%  "

$ws.Range("D9").Value = "'44.79"
$ws.Range("E9").Value = "  -6.59%  "

$ws.Range("E10").Value = "  -3.72%  "

$ws.Range("D11").Value = "'0.07219"
$ws.Range("E11").Value = "  -3.97%  "

$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  -1.00%  "

$ws.Range("D13").Value = "'22.36"
$ws.Range("E13").Value = "  +0.14%  "

$ws.Range("D14").Value = "'6.140"
$ws.Range("E14").Value = "  -4.59%  "

$ws.Range("E15").Value = "  +0.04%  "

$ws.Range("D16").Value = "1.758.76"
$ws.Range("E16").Value = "  -1.66%  "

$ws.Range("E17").Value = "  -2.98%  "

$ws.Range("D18").Value = "'0.06607"
$ws.Range("E18").Value = "  -1.70%  "

$ws.Range("E19").Value = "  -5.36%  "

$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.47%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'6.224"
$ws.Range("E21").Value = "  -4.60%  "

$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").Value = "'16.63"
$ws.Range("E22").Value = "  -5.91%  "

$ws.Range("D23").Value = "27.973.24"
$ws.Range("E23").Value = "  +1.13%  "

$ws.Range("E24").Value = "  -5.80%  "

$ws.Range("D25").Value = "'2.387"
$ws.Range("E25").Value = "  -0.96%  "

$ws.Range("D26").Value = "'152.30"
$ws.Range("E26").Value = "  -1.72%  "

$ws.Range("D27").Value = "'19.78"
$ws.Range("E27").Value = "  -6.15%  "

$ws.Range("E28").Value = "  -7.46%  "

$ws.Range("D29").Value = "1.959.26"
$ws.Range("E29").Value = "  -1.62%  "

$ws.Range("D30").Value = "'1.268"
$ws.Range("E30").Value = "  -15.42%  "

$ws.Range("D31").Value = "'131.73"
$ws.Range("E31").Value = "  -3.88%  "

$ws.Range("D32").Value = "'4.016"
$ws.Range("E32").Value = "  +0.20%  "

$ws.Range("D33").Value = "'5.840"
$ws.Range("E33").Value = "  -6.77%  "

$ws.Range("D34").Value = "'0.08813"
$ws.Range("E34").Value = "  -0.23%  "

$ws.Range("D35").Value = "'12.20"
$ws.Range("E35").Value = "  -6.15%  "

$ws.Range("D36").Value = "'0.6616"
$ws.Range("E36").Value = "  -3.87%  "

$ws.Range("E37").Value = "  -3.80%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02296"
$ws.Range("E38").Value = "  -6.38%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'5.167"
$ws.Range("E39").Value = "  -5.51%  "

$ws.Range("D40").Value = "'0.2111"
$ws.Range("E40").Value = "  -5.19%  "

$ws.Range("D41").Value = "'1.500"
$ws.Range("E41").Value = "  -5.51%  "

$ws.Range("D42").Value = "'1.209"
$ws.Range("E42").Value = "  -4.51%  "

$ws.Range("D43").Value = "'7.954"
$ws.Range("E43").Value = "  -6.37%  "

$ws.Range("D44").Value = "'0.9997"
$ws.Range("E44").Value = "  -0.66%  "

$ws.Range("D45").Value = "'13.80"
$ws.Range("E45").Value = "  -6.05%  "

$ws.Range("D46").Value = "'3.823"
$ws.Range("E46").Value = "  -1.24%  "

$ws.Range("D47").Value = "'0.6040"
$ws.Range("E47").Value = "  -4.84%  "

$ws.Range("D48").Value = "'126.17"
$ws.Range("E48").Value = "  -5.51%  "

$ws.Range("D49").Value = "'2.006"
$ws.Range("E49").Value = "  -6.35%  "

$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "'1.177"
$ws.Range("E50").Value = "  +1.20%  "

$ws.Range("B51").Value = "Flow"
$ws.Range("C51").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D51").Value = "'1.113"
$ws.Range("E51").Value = "  +3.10%  "

